# Updated BackLog for Sprint 2
# Adds "Sprint" (col E) and "Priority" (col F) data to the Product Backlog,
# renames story #4 from "Select 3 Pairs of Dice" to "Pair Dice", and moves
# the active selection to C23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Story name fix (row 6 / Story ID 4): "Select 3 Pairs of Dice" -> "Pair Dice"
$ws.Range("B6").Value = "Pair Dice"

# --- Sprint (E) / Priority (F) values for every backlog item ---
# Rows 4-7 already have blank, formatted E/F cells; typing into them keeps
# their existing look. Rows 8-20 get their first-ever E/F entries, so copy
# the formatting from a neighboring already-styled cell onto them before
# writing the values, matching how the sheet reads after manual editing.

$fmtSrc = $ws.Range("E4")

foreach ($row in 8..20) {
    $eCell = $ws.Cells.Item($row, 5)
    $fCell = $ws.Cells.Item($row, 6)

    $fmtSrc.Copy()
    $eCell.PasteSpecial(-4122)
    $fmtSrc.Copy()
    $fCell.PasteSpecial(-4122)
}

$sprintData = @{
    3  = 2
    4  = 2
    5  = 3
    6  = 2
    7  = 2
    8  = 3
    9  = 3
    10 = 3
    11 = 3
    12 = 3
    13 = 3
    14 = 2
    15 = 2
    16 = 3
    17 = 3
    18 = 2
    19 = 2
    20 = 2
}

# Priority text is written in "High" -> "Medium" -> "Low" first-seen order
# (row 4, then row 5, then row 3) so the shared-string table is built up in
# that order, then the rest of the rows follow in natural order.
$priorityOrder = @(4, 5, 3, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)
$priorityData = @{
    3  = "Low"
    4  = "High"
    5  = "Medium"
    6  = "High"
    7  = "High"
    8  = "Medium"
    9  = "Medium"
    10 = "Medium"
    11 = "Medium"
    12 = "Medium"
    13 = "Medium"
    14 = "High"
    15 = "High"
    16 = "Low"
    17 = "High"
    18 = "High"
    19 = "High"
    20 = "Low"
}

foreach ($row in $priorityOrder) {
    $ws.Cells.Item($row, 6).Value = $priorityData[$row]
}

foreach ($row in 3..20) {
    $ws.Cells.Item($row, 5).Value = $sprintData[$row]
}

# --- Active cell / selection moves to C23 ---
$ws.Range("C23").Select()
